# Fix a typo in the "Categóricos" sheet: the header of the VISTAS table
# (cell A2, which is also the Table1/"VISTAS" column-1 header) was
# misspelled "cvevisita" and is corrected to "cvevista".
#
# Editing the table header cell through the object model also re-syncs
# the underlying Excel Table's column name automatically, matching the
# corresponding <tableColumn name="..."/> update in xl/tables/table1.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categóricos")

$ws.Range("A2").Value = "cvevista"
